$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (E..T) mapped to indexes 5..20
# Row => hashtable of column letter => new value
$updates = @{
    2  = @{ E=3; G=46.926932; H=140.780796; I=0.3970508522302297; J=0.3970508522302297; K=3; M=20.7608; N=62.2824; O=0.3489417963420857; P=0.3489417963420857; Q=974.2406498655999; R=8768.165848790401; S=0.1385476376163724; T=0.1385476376163724 }
    3  = @{ E=3; G=46.926932; H=140.780796; I=0.3970508522302297; J=0.3970508522302297; K=3; M=23.032996; N=69.09898799999999; O=0.3871322395755498; P=0.3871322395755498; Q=1080.867837048272; R=9727.810533434447; S=0.1537111856492695; T=0.1537111856492695 }
    4  = @{ E=3; G=46.926932; H=140.780796; I=0.3970508522302297; J=0.3970508522302297; K=3; M=15.70265933333333; N=47.107978; O=0.2639259640823645; P=0.2639259640823645; Q=736.8776267544988; R=6631.898640790489; S=0.1047920289645878; T=0.1047920289645878 }
    5  = @{ E=3; G=30.21862233333333; H=90.655867; I=0.2556811033517694; J=0.2556811033517694; K=3; M=20.7608; N=62.2824; O=0.3489417963420857; P=0.3489417963420857; Q=627.3627745378666; R=5646.2649708408; S=0.08921782349429287; T=0.08921782349429287 }
    6  = @{ E=3; G=30.21862233333333; H=90.655867; I=0.2556811033517694; J=0.2556811033517694; K=3; M=23.032996; N=69.09898799999999; O=0.3871322395755498; P=0.3871322395755498; Q=696.0254073291773; R=6264.228665962595; S=0.09898239815771807; T=0.09898239815771809 }
    7  = @{ E=3; G=30.21862233333333; H=90.655867; I=0.2556811033517694; J=0.2556811033517694; K=3; M=15.70265933333333; N=47.107978; O=0.2639259640823645; P=0.2639259640823645; Q=474.5127320229918; R=4270.614588206926; S=0.0674808816997584; T=0.0674808816997584 }
    8  = @{ E=3; G=41.04316566666667; H=123.129497; I=0.3472680444180009; J=0.3472680444180009; K=3; M=20.7608; N=62.2824; O=0.3489417963420857; P=0.3489417963420857; Q=852.0889537725335; R=7668.8005839528; S=0.1211763352314204; T=0.1211763352314204 }
    9  = @{ E=3; G=41.04316566666667; H=123.129497; I=0.3472680444180009; J=0.3472680444180009; K=3; M=23.032996; N=69.09898799999999; O=0.3871322395755498; P=0.3871322395755498; Q=945.3470706276707; R=8508.123635649035; S=0.1344386557685622; T=0.1344386557685622 }
    10 = @{ E=3; G=41.04316566666667; H=123.129497; I=0.3472680444180009; J=0.3472680444180009; K=3; M=15.70265933333333; N=47.107978; O=0.2639259640823645; P=0.2639259640823645; Q=644.4868484252297; R=5800.381635827067; S=0.09165305341801824; T=0.09165305341801824 }
}

foreach ($rowNum in $updates.Keys) {
    $rowVals = $updates[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $cellRef = "$col$rowNum"
        $ws.Range($cellRef).Value = $rowVals[$col]
    }
}
